$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 817
$ws.Range("F3").Value = 14714
$ws.Range("F4").Value = 14945
$ws.Range("F5").Value = 6046
$ws.Range("F8").Value = 1638
$ws.Range("F9").Value = 494
$ws.Range("F11").Value = 1294
$ws.Range("F12").Value = 1965
$ws.Range("F13").Value = 945
$ws.Range("F14").Value = 43
$ws.Range("F15").Value = 2342
$ws.Range("F16").Value = 611
$ws.Range("F18").Value = 3619
$ws.Range("F20").Value = 345
$ws.Range("F21").Value = 2665
$ws.Range("F22").Value = 678
$ws.Range("F24").Value = 1343
$ws.Range("F25").Value = 1913
$ws.Range("F26").Value = 1135
$ws.Range("F27").Value = 1608
$ws.Range("F28").Value = 344
$ws.Range("F29").Value = 173
$ws.Range("F30").Value = 7495
$ws.Range("F31").Value = 5194
$ws.Range("F32").Value = 330
$ws.Range("F34").Value = 724
$ws.Range("F35").Value = 724
$ws.Range("F36").Value = 3396
$ws.Range("F39").Value = 359
$ws.Range("F40").Value = 156
$ws.Range("F41").Value = 116
$ws.Range("F42").Value = 4507
$ws.Range("F43").Value = 723
$ws.Range("F44").Value = 29
$ws.Range("F45").Value = 346

$ws = $wb.Worksheets.Item("演出")
$ws.Range("G3").Value = "不可售"
$ws.Range("F13").Value = 23
$ws.Range("F20").Value = 114

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 8008
$ws.Range("F4").Value = 1116

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 8008
$ws.Range("F3").Value = 817
$ws.Range("F5").Value = 1116
$ws.Range("F6").Value = 14714
$ws.Range("F7").Value = 14714
$ws.Range("F8").Value = 6046
$ws.Range("F12").Value = 1638
$ws.Range("F13").Value = 494
$ws.Range("F14").Value = 1294
$ws.Range("F15").Value = 1965
$ws.Range("F17").Value = 43
$ws.Range("F19").Value = 3619
$ws.Range("F20").Value = 345
$ws.Range("F21").Value = 678
$ws.Range("F24").Value = 1913
$ws.Range("F26").Value = 23
$ws.Range("F30").Value = 1608
$ws.Range("F32").Value = 344
$ws.Range("F33").Value = 173
$ws.Range("F34").Value = 7495
$ws.Range("F35").Value = 5194
$ws.Range("F36").Value = 330
$ws.Range("F37").Value = 724
$ws.Range("F38").Value = 3396
$ws.Range("F41").Value = 359
$ws.Range("F43").Value = 116
$ws.Range("F44").Value = 4507
$ws.Range("F45").Value = 723
$ws.Range("F46").Value = 29
$ws.Range("F47").Value = 346
